$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''317.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''3.94%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''39.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''2.37%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''5.149'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''0.69%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.08231'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''1.85%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''2.074'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''7.85%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''8.358'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''4.45%'
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9399'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''0.97%'
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.1364'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''-6.39%'
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1988'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''3.29%'
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.09153'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''0.39%'
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03509'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''-0.02%'
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09813'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''0.15%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001406'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''0.81%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006042'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''1.86%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.691'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''-2.29%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''4.337'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''3.48%'
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''3.336'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''-3.33%'
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''0.3478'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''0.58%'
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''0.1315'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''-2.52%'
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''4.953'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''5.88%'
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.2450'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''1.36%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.04359'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''-0.62%'
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''-0.86%'
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''13.04%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.0001298'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''-0.42%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.0003998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''-10.10%'
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = '''0.02234'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''9.88%'
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.05222'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = '''0.007762'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''2.70%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.009694'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''-0.50%'
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''4.50%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.002046'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''-3.78%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.009660'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''-2.78%'
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.00006625'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''7.07%'
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = '''-0.40%'
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''1.73%'
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.001689'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''-6.42%'
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''-0.40%'
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''-0.40%'
$ws.Range("E51").Style = "Normal"
